# Add a new row of expense data (category labels) below the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "chik"
$ws.Cells.Item(2, 2).Value = "Food"

# Price and Date look numeric/date-like to Excel's auto-detection, so force
# them to be entered as plain text (matching the source data) and then
# restore the default cell style so no extra formatting sticks around.
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "$23.67"
$ws.Cells.Item(2, 3).Style = "Normal"

$ws.Cells.Item(2, 4).Value = "High"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "6/10/24"
$ws.Cells.Item(2, 5).Style = "Normal"
